$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D ("Koronawirus" - test result), shifting the
# existing "posiadane pieniądze" column (and its data) from D to E.
$ws.Range("D1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Koronawirus"

# New "Koronawirus" (test result) column values
$ws.Range("D2").Value = "brak"
$ws.Range("D3").Value = "brak"
$ws.Range("D4").Value = "pozytywny"
$ws.Range("D5").Value = "brak"

# Updated "posiadane pieniądze" values (now shifted to column E)
$ws.Range("E2").Value = 23
$ws.Range("E3").Value = 220
$ws.Range("E4").Value = 100
$ws.Range("E5").Value = 23
